$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new values look numeric,
# so Excel keeps them as text (matching the original inline-string data)
# instead of silently converting them into floating point numbers.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply the updated cell values row by row.
# Row 2
$ws.Range("D2").Value = '63.004.71'
$ws.Range("E2").Value = '  +5.10%  '

# Row 3
$ws.Range("D3").Value = '3.368.27'
$ws.Range("E3").Value = '  +5.77%  '

# Row 4
$ws.Range("E4").Value = '  +0.03%  '

# Row 5
$ws.Range("D5").Value = '572.78'
$ws.Range("E5").Value = '  +7.19%  '

# Row 6
$ws.Range("D6").Value = '152.74'
$ws.Range("E6").Value = '  +5.19%  '

# Row 8
$ws.Range("D8").Value = '3.370.54'
$ws.Range("E8").Value = '  +5.63%  '

# Row 9
$ws.Range("D9").Value = '0.525'
$ws.Range("E9").Value = '  -0.26%  '

# Row 10
$ws.Range("D10").Value = '7.43'
$ws.Range("E10").Value = '  +1.67%  '

# Row 11
$ws.Range("D11").Value = '0.118'
$ws.Range("E11").Value = '  +5.60%  '

# Row 12
$ws.Range("D12").Value = '0.435'
$ws.Range("E12").Value = '  +1.31%  '

# Row 13
$ws.Range("D13").Value = '3.951.65'
$ws.Range("E13").Value = '  +5.83%  '

# Row 14
$ws.Range("E14").Value = '  +0.26%  '

# Row 15
$ws.Range("B15").Value = 'ShibaInu'
$ws.Range("C15").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D15").Value = '0.0000180'
$ws.Range("E15").Value = '  +4.54%  '

# Row 16
$ws.Range("B16").Value = 'Avalanche'
$ws.Range("C16").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D16").Value = '26.93'
$ws.Range("E16").Value = '  +4.17%  '

# Row 17
$ws.Range("D17").Value = '63.098.11'
$ws.Range("E17").Value = '  +5.19%  '

# Row 18
$ws.Range("D18").Value = '3.331.72'
$ws.Range("E18").Value = '  +4.60%  '

# Row 19
$ws.Range("D19").Value = '6.31'
$ws.Range("E19").Value = '  +0.89%  '

# Row 20
$ws.Range("D20").Value = '13.89'
$ws.Range("E20").Value = '  +5.05%  '

# Row 21
$ws.Range("D21").Value = '8.39'
$ws.Range("E21").Value = '  +2.53%  '

# Row 22
$ws.Range("D22").Value = '384.92'
$ws.Range("E22").Value = '  +4.63%  '

# Row 23
$ws.Range("D23").Value = '1.00'
$ws.Range("E23").Value = '  +0.05%  '

# Row 24
$ws.Range("E24").Value = '  +2.46%  '

# Row 25
$ws.Range("E25").Value = '  +1.30%  '

# Row 26
$ws.Range("E26").Value = '  +6.61%  '

# Row 27
$ws.Range("D27").Value = '9.26'
$ws.Range("E27").Value = '  +6.72%  '

# Row 28
$ws.Range("D28").Value = '0.0₃0967'
$ws.Range("E28").Value = '  +11.83%  '

# Row 29
$ws.Range("D29").Value = '0.999'
$ws.Range("E29").Value = '  +0.06%  '

# Row 30
$ws.Range("E30").Value = '  +7.26%  '

# Row 31
$ws.Range("B31").Value = 'RenderToken'
$ws.Range("C31").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D31").Value = '6.39'
$ws.Range("E31").Value = '  +5.81%  '

# Row 32
$ws.Range("B32").Value = 'EthereumClassic'
$ws.Range("C32").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D32").Value = '23.02'
$ws.Range("E32").Value = '  +3.21%  '

# Row 33
$ws.Range("B33").Value = 'Fetch.AI'
$ws.Range("C33").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D33").Value = '1.31'
$ws.Range("E33").Value = '  +10.55%  '

# Row 34
$ws.Range("B34").Value = 'NEARProtocol'
$ws.Range("C34").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D34").Value = '5.57'
$ws.Range("E34").Value = '  +5.59%  '

# Row 35
$ws.Range("D35").Value = '6.69'
$ws.Range("E35").Value = '  +2.07%  '

# Row 36
$ws.Range("E36").Value = '  +9.92%  '

# Row 37
$ws.Range("D37").Value = '158.31'
$ws.Range("E37").Value = '  +1.41%  '

# Row 38
$ws.Range("E38").Value = '  +12.42%  '

# Row 39
$ws.Range("D39").Value = '27.31'
$ws.Range("E39").Value = '  +4.66%  '

# Row 40
$ws.Range("D40").Value = '2.894.72'
$ws.Range("E40").Value = '  +2.91%  '

# Row 41
$ws.Range("D41").Value = '0.0327'
$ws.Range("E41").Value = '  +10.06%  '

# Row 42
$ws.Range("D42").Value = '0.0741'
$ws.Range("E42").Value = '  +5.42%  '

# Row 43
$ws.Range("D43").Value = '40.77'
$ws.Range("E43").Value = '  +2.81%  '

# Row 44
$ws.Range("E44").Value = '  +4.38%  '

# Row 45
$ws.Range("D45").Value = '4.22'
$ws.Range("E45").Value = '  +0.07%  '

# Row 46
$ws.Range("B46").Value = 'RenzoRestakedETH'
$ws.Range("C46").Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range("D46").Value = '3.418.08'
$ws.Range("E46").Value = '  +5.90%  '

# Row 47
$ws.Range("B47").Value = 'ONDO'
$ws.Range("C47").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D47").Value = '1.04'
$ws.Range("E47").Value = '  +5.84%  '

# Row 48
$ws.Range("B48").Value = 'Bittensor'
$ws.Range("C48").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D48").Value = '301.99'
$ws.Range("E48").Value = '  +14.68%  '

# Row 49
$ws.Range("B49").Value = 'InjectiveProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D49").Value = '21.89'
$ws.Range("E49").Value = '  +6.08%  '

# Row 50
$ws.Range("E50").Value = '  -1.98%  '

# Row 51
$ws.Range("D51").Value = '6.29'
$ws.Range("E51").Value = '  +2.31%  '
